# Update the dSF column (F) with repulled/recalculated data values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value  = 1
$ws.Range("F3").Value  = -2
$ws.Range("F4").Value  = -1
$ws.Range("F5").Value  = 3
$ws.Range("F6").Value  = -3
$ws.Range("F7").Value  = 9
$ws.Range("F8").Value  = 1
$ws.Range("F9").Value  = 2
$ws.Range("F11").Value = 1
$ws.Range("F12").Value = 5
$ws.Range("F13").Value = -1
$ws.Range("F14").Value = 2

$wb.Save()
